$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.195.60"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.825.54"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'235.77"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "'0.6099"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.07086"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").Value = "'0.2804"
$ws.Range("E9").Value = "  -3.01%  "
$ws.Range("D10").Value = "'23.51"
$ws.Range("E10").Value = "  -6.03%  "
$ws.Range("D11").Value = "'0.07660"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "1.824.14"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "'4.806"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "'0.000009996"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "'0.6310"
$ws.Range("E15").Value = "  -6.54%  "
$ws.Range("D16").Value = "2.068.65"
$ws.Range("E16").Value = "  -1.04%  "
$ws.Range("D17").Value = "'78.59"
$ws.Range("E17").Value = "  -3.66%  "
$ws.Range("D18").Value = "'5.858"
$ws.Range("E18").Value = "  -6.02%  "
$ws.Range("D19").Value = "29.228.92"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("D20").Value = "'226.85"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "'11.76"
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("D23").Value = "'7.000"
$ws.Range("E23").Value = "  -4.61%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'155.42"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'8.034"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("D27").Value = "'0.1307"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").Value = "'16.60"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").Value = "'1.490"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").Value = "'0.06346"
$ws.Range("E30").Value = "  -12.44%  "
$ws.Range("D31").Value = "'1.448"
$ws.Range("E31").Value = "  -2.12%  "
$ws.Range("D32").Value = "'3.815"
$ws.Range("E32").Value = "  -5.53%  "
$ws.Range("D33").Value = "'3.787"
$ws.Range("E33").Value = "  -6.11%  "
$ws.Range("D34").Value = "'1.125"
$ws.Range("E34").Value = "  -1.21%  "
$ws.Range("D35").Value = "'1.741"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").Value = "'0.6452"
$ws.Range("E36").Value = "  -7.15%  "
$ws.Range("D37").Value = "'2.548"
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'2.738"
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.212.51"
$ws.Range("E39").Value = "  -1.75%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01738"
$ws.Range("E40").Value = "  -5.46%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.533"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").Value = "'0.9133"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.980.03"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'100.66"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'62.57"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "'0.00000000115"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'8.557"
$ws.Range("E48").Value = "  -3.94%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.604"
$ws.Range("E49").Value = "  -6.04%  "
$ws.Range("D50").Value = "'0.4564"
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'0.05520"
$ws.Range("E51").Value = "  -2.55%  "
